# McuOpenPnP Bill of Material: "updated number of idlers"
# Row 40 is the "GT2 Idler Pulley 20T 5mm_Toothed" line item.
# The idler quantity went from 4 to 6, and a new note was added
# explaining the quantity depends on configuration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the idler pulley quantity (column C) from 4 to 6.
# The E column (Costs) recalculates automatically via its shared formula,
# as does the E2 total (SUM(E4:E66)).
$ws.Range("C40").Value = 6

# Add a note in column F (Comment) explaining the quantity.
$ws.Range("F40").Value = "Depends on your configuration. You need at least 3."

# Reflect the author's last active selection in the saved view.
$ws.Range("C37").Select()
